# ChargesGrid update: charge_button_dict -> set_plea
# Shifts the OVI charge into column 2 as a new breath-test charge, shifts
# "Traffic Control Lights" into column 3 (replacing the dismissed marijuana
# charge), zeroes out the fine amounts/suspensions tied to the old column 1/3
# charges, and bumps the arraignment / paid-in-full dates by one day.

$d = $word.ActiveDocument

# Find.Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards, MatchSoundsLike,
#              MatchAllWordForms, Forward, Wrap, Format, ReplaceWith, Replace)
# Replace: 0=wdReplaceNone, 1=wdReplaceOne, 2=wdReplaceAll
# wdReplaceOne is used throughout so that the duplicated "MM" cell text only
# touches the intended (first / column-2) occurrence.

# 1. Arraignment date sentence.
$d.Content.Find.Execute(" on March 19, 2022.", $true, $false, $false, $false, $false, `
                         $true, 1, $false, " on March 20, 2022.", 1) | Out-Null

# 2. Offense row, column 2 (was "Traffic Control Lights").
$d.Content.Find.Execute("Traffic Control Lights", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "OVI Breath 1st .17  Above", 1) | Out-Null

# 3. Offense row, column 3 (was "Possession of Marijuana less than 100 grams").
$d.Content.Find.Execute("Possession of Marijuana less than 100 grams", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Traffic Control Lights", 1) | Out-Null

# 4. Statute/Ord. row, column 2 (was "4511.13C").
$d.Content.Find.Execute("4511.13C", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "4511.19A1H*", 1) | Out-Null

# 5. Statute/Ord. row, column 3 (was "2925.11(C)(3)(a)").
$d.Content.Find.Execute("2925.11(C)(3)(a)", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "4511.13C", 1) | Out-Null

# 6. Degree row, column 2 (first "MM" in the table, was "MM").
$d.Content.Find.Execute("MM", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "M1", 1) | Out-Null

# 7. Fine Amount row, column 1 (was "$ 34").
$d.Content.Find.Execute("$ 34", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "$ 0", 1) | Out-Null

# 8. Fine Amount row, column 3 (was "$ 123").
$d.Content.Find.Execute("$ 123", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "$ 0", 1) | Out-Null

# 9. Fines Suspended row, column 1 (was "$ 3").
$d.Content.Find.Execute("$ 3", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "$ 0", 1) | Out-Null

# 10. Fines Suspended row, column 3 (was "$ 1").
$d.Content.Find.Execute("$ 1", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "$ 0", 1) | Out-Null

# 11. "Paid in full by" date.
$d.Content.Find.Execute("March 19, 2022", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "March 20, 2022", 1) | Out-Null
